$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers; force text format
# so Excel does not silently convert them to numeric cells (they were
# stored as text/inlineStr strings in the original workbook).
$textForceAddrs = @("D5", "D6", "D8", "D9", "D13", "D18", "D19", "D20", "D23", "D24", "D26", "D30", "D31", "D37", "D41", "D45", "D46", "D48")
foreach ($addr in $textForceAddrs) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated numeric-looking (but text) values
$ws.Range("D5").Value = "597.43"
$ws.Range("D6").Value = "145.17"
$ws.Range("D8").Value = "0.502"
$ws.Range("D9").Value = "6.98"
$ws.Range("D13").Value = "33.46"
$ws.Range("D18").Value = "6.69"
$ws.Range("D19").Value = "431.50"
$ws.Range("D20").Value = "13.48"
$ws.Range("D23").Value = "81.76"
$ws.Range("D24").Value = "10.85"
$ws.Range("D26").Value = "11.73"
$ws.Range("D30").Value = "6.91"
$ws.Range("D31").Value = "26.59"
$ws.Range("D37").Value = "2.98"
$ws.Range("D41").Value = "42.13"
$ws.Range("D45").Value = "133.66"
$ws.Range("D46").Value = "361.91"
$ws.Range("D48").Value = "23.59"

# Apply remaining text value updates
$ws.Range("D2").Value = "61.381.49"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.928.61"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "3.413.48"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "61.356.12"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.925.80"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "0.0₃0884"
$ws.Range("E34").Value = "  +3.36%  "
$ws.Range("E35").Value = "  -0.46%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  +5.62%  "
$ws.Range("E42").Value = "  -1.99%  "
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "2.699.77"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  +2.31%  "
$ws.Range("E46").Value = "  -3.79%  "
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("E51").Value = "  -2.27%  "
